# Word session / active document (pre-seeded by the host).
$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Helper: locate a paragraph (1-based index, at or after $after) whose text
# contains $needle. Uses .Contains() (case-sensitive) rather than -like,
# since e.g. "Event_date" / "event_date" both appear in this document and
# must not be confused with one another.
function Find-ParaIndex($doc, [string]$needle, [int]$after = 0) {
    for ($i = $after + 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# "User:" section — First_name / Last_name paragraphs get spell-checker
# proofErr wrapping around the flagged (underscored) words, and each
# paragraph's single run is split into several runs.
# ---------------------------------------------------------------------------

$firstNameIdx = Find-ParaIndex $d "First_name"
$firstNameXml = '<w:p ' + $wNs + '>' + `
    '<w:r><w:tab/></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>First_name</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve">: no numbers, no special chars, at least 3 chars, </w:t></w:r>' + `
  '</w:p>'
$null = $d.Paragraphs($firstNameIdx).Range.InsertXML($firstNameXml)

$lastNameIdx = Find-ParaIndex $d "_name: no numbers" $firstNameIdx
$lastNameXml = '<w:p ' + $wNs + '>' + `
    '<w:r><w:tab/></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Last</w:t></w:r>' + `
    '<w:r><w:t>_name</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve">: no numbers, no special chars, at least 3 chars, </w:t></w:r>' + `
  '</w:p>'
$null = $d.Paragraphs($lastNameIdx).Range.InsertXML($lastNameXml)

# ---------------------------------------------------------------------------
# "Event:" section — Event_date / Event_time paragraphs get the same
# proofErr treatment on every underscored token.
# ---------------------------------------------------------------------------

$eventDateIdx = Find-ParaIndex $d "Event_date" $lastNameIdx
$eventDateXml = '<w:p ' + $wNs + '>' + `
    '<w:r><w:tab/></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Event_date</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve">: present, </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>date_format</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
  '</w:p>'
$null = $d.Paragraphs($eventDateIdx).Range.InsertXML($eventDateXml)

$eventTimeIdx = Find-ParaIndex $d "Event_time" $eventDateIdx
$eventTimeXml = '<w:p ' + $wNs + '>' + `
    '<w:r><w:tab/></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Event_time</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve">: present, </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>time_format</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
  '</w:p>'
$null = $d.Paragraphs($eventTimeIdx).Range.InsertXML($eventTimeXml)

# ---------------------------------------------------------------------------
# Trailing "Added validation Event model:" / event_date / event_time
# paragraphs are removed entirely, collapsing down to a single empty
# paragraph that just carries the _GoBack bookmark.
# ---------------------------------------------------------------------------

$headerIdx = Find-ParaIndex $d "Added validation Event model" $eventTimeIdx
$lastTailIdx = Find-ParaIndex $d "the current time" $headerIdx

$start = $d.Paragraphs($headerIdx).Range.Start
$end = $d.Paragraphs($lastTailIdx).Range.End
$tailRange = $d.Range($start, $end)

$tailXml = '<w:p ' + $wNs + '>' + `
    '<w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
  '</w:p>'
$null = $tailRange.InsertXML($tailXml)

Write-Output "done"
